$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compartments")

# Insert a new column before column E (shifts Database references/Comments/Created right)
$ws.Columns.Item(5).Insert()

# Rename the "Initial volume" header to "Mean volume"
$ws.Range("D1").Value = "Mean volume"

# New column header + unit values
$ws.Range("E1").Value = "Mean volume units"
$ws.Range("E2").Value = "L"
$ws.Range("E3").Value = "L"

# Re-apply the autofilter so it covers the newly inserted column
$ws.AutoFilterMode = $false
$ws.Range("A1:G3").AutoFilter()

# Make Compartments the active/selected sheet and cell
$ws.Select()
$ws.Range("E2").Select()
